$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The song list (A2:D52) is sorted alphabetically (A -> Z) by the
# "Song Title" column (B), rows moving together as whole records.
# Header parameter is explicitly xlNo (2) since row 1 is a separate
# header row and is not part of the sorted range.
$rng = $ws.Range("A2:D52")
$key = $ws.Range("B2:B52")

$rng.Sort($key, 1, $null, $null, 1, $null, 1, 2)
